# "updated docs for pages"
#
# On the "Metadata" sheet:
#   - B7 (the "Experimental" row) gets the literal text value "true"
#     (it was previously blank).
#   - B8 (the "Date" row) is updated to the new generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Plain `Range.Value = "true"` would be auto-recognized by Excel as the
# boolean TRUE instead of the text "true". To force a genuine text value
# (matching the shared-string cell the workbook expects) we type it with a
# leading apostrophe into a scratch cell first - exactly like typing 'true
# into a cell in the Excel UI - then copy only the *value* over to B7 so
# B7 keeps its original cell formatting/style untouched.
$scratch = $ws.Range("Z1")
$scratch.Value = "'true"
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)  # xlPasteValues
$scratch.Clear()

$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
